# "aded Candidate Sheet + ordered events" — add the missing "Fill Out
# Candidate Information Sheet" confirmation-event row to the
# "Confirmation Events" sheet (row 5 was reserved/blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmation Events")

# Seed A5/B5/C5 with the same cell formatting already used by the other
# event rows (name/instructions columns share the plain "s=2" style, the
# due_date column uses the date-formatted "s=4" style) before writing the
# values, so the new row matches its siblings exactly instead of Excel
# minting a brand new ad-hoc style for the date cell.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("C4").Copy($ws.Range("C5"))

$ws.Range("A5").Value = "Fill Out Candidate Information Sheet"
$ws.Range("B5").Value = "2/16/2016"
$ws.Range("C5").Value = "<p><em><strong>simple text</strong></em></p>"
